$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 1
    6  = 2
    7  = 0
    8  = 1
    9  = 2
    10 = 2
    11 = 1
    12 = 0
    13 = 2
    14 = 2
    15 = 1
    16 = 2
    17 = 2
    18 = 0
    19 = 2
    20 = 0
    21 = 2
    22 = 0
    23 = 1
    24 = 1
    25 = 2
    26 = 1
    27 = 2
    28 = 3
    29 = 3
    30 = 0
    31 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
